# Apply the benchmark-table numeric updates described by the commit:
# "Fixed README.md stats and docx preparation for all Renaissance -
#  JDK 17 - Z GC tests"
#
# The document is a single-column table; each change below targets one
# row (cell) by its 1-based row index and overwrites that cell's text
# in place, which preserves the existing run formatting
# (Times New Roman, sz 22).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "208"
    6  = "0.00015"
    7  = "0.00005"
    8  = "0.00003"
    9  = "0.00004"
    10 = "0.00005"
    11 = "0.00012"
    12 = "0.01133"
    44 = "99.99"
    45 = "0.01"
    46 = "174"
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $updates[$row]
}
